$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "103÷7=14, 5";   New = "938÷4=234, 2" },
    @{ Old = "401÷9=44, 5";   New = "948÷4=237, 0" },
    @{ Old = "454÷6=75, 4";   New = "318÷4=79, 2" },
    @{ Old = "219÷4=54, 3";   New = "910÷9=101, 1" },
    @{ Old = "539÷8=67, 3";   New = "298÷3=99, 1" },
    @{ Old = "588÷9=65, 3";   New = "861÷4=215, 1" },
    @{ Old = "991÷4=247, 3";  New = "853÷4=213, 1" },
    @{ Old = "892÷7=127, 3";  New = "500÷9=55, 5" },
    @{ Old = "702÷5=140, 2";  New = "817÷5=163, 2" },
    @{ Old = "712÷8=89, 0";   New = "927÷2=463, 1" },
    @{ Old = "260÷3=86, 2";   New = "147÷8=18, 3" },
    @{ Old = "377÷2=188, 1";  New = "375÷3=125, 0" },
    @{ Old = "299÷5=59, 4";   New = "260÷2=130, 0" },
    @{ Old = "867÷4=216, 3";  New = "608÷4=152, 0" },
    @{ Old = "790÷2=395, 0";  New = "227÷7=32, 3" },
    @{ Old = "296÷2=148, 0";  New = "545÷6=90, 5" },
    @{ Old = "658÷8=82, 2";   New = "644÷8=80, 4" },
    @{ Old = "963÷6=160, 3";  New = "773÷6=128, 5" },
    @{ Old = "770÷7=110, 0";  New = "602÷2=301, 0" },
    @{ Old = "869÷4=217, 1";  New = "784÷4=196, 0" },
    @{ Old = "684÷4=171, 0";  New = "907÷5=181, 2" },
    @{ Old = "692÷4=173, 0";  New = "734÷4=183, 2" },
    @{ Old = "664÷6=110, 4";  New = "491÷3=163, 2" },
    @{ Old = "713÷4=178, 1";  New = "561÷2=280, 1" },
    @{ Old = "934÷6=155, 4";  New = "441÷3=147, 0" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
